$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A70:A71").NumberFormat = "@"

$ws.Range("A70").Value = "2024-05-13"
$ws.Range("B70").Value = "12:44:57"
$ws.Range("C70").Value = "Fallo tornillo"
$ws.Range("D70").Value = "-"
$ws.Range("E70").Value = "-"
$ws.Range("F70").Value = "-"
$ws.Range("G70").Value = "-"

$ws.Range("A71").Value = "2024-05-13"
$ws.Range("B71").Value = "12:45:02"
$ws.Range("C71").Value = "-"
$ws.Range("D71").Value = "Tornillo atascado en tolva"
$ws.Range("E71").Value = "-"
$ws.Range("F71").Value = "-"
$ws.Range("G71").Value = "-"
